# Applies updated bus voltage magnitude (vm_pu) results for the 380 kV case
# (commit: "case with 380 kV done") to Sheet1, rows 2-25, columns B-F and I-N.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"=1.02; "C"=1.022515186447222; "D"=1.027322072015636; "E"=1.047234263971629; "F"=1.051171857721807; "I"=1.028889589830776; "J"=1.027700590172546; "K"=1.030142139068718; "L"=1.049997471588898; "M"=1.05392409815486; "N"=1.029160042889497 }
    3 = @{ "B"=1.02; "C"=1.023519491162126; "D"=1.028068211800912; "E"=1.048536183682562; "F"=1.052558003099678; "I"=1.02906534233552; "J"=1.028342692500329; "K"=1.030696146038918; "L"=1.05110992339183; "M"=1.055121361341492; "N"=1.029803057076235 }
    4 = @{ "B"=1.02; "C"=1.024168644843595; "D"=1.02854985966099; "E"=1.049379044747758; "F"=1.053455258724909; "I"=1.029176753343894; "J"=1.02875691907697; "K"=1.031052812161136; "L"=1.051829639490132; "M"=1.05589586779679; "N"=1.030217871902127 }
    5 = @{ "B"=1.02; "C"=1.024441382265813; "D"=1.028752067511025; "E"=1.04973348934201; "F"=1.053832545763459; "I"=1.029223036488132; "J"=1.028930759126268; "K"=1.031202320126545; "L"=1.052132182791583; "M"=1.056221423853749; "N"=1.030391958824235 }
    6 = @{ "B"=1.02; "C"=1.024487166312043; "D"=1.028786002844958; "E"=1.0497930084228; "F"=1.053895898727124; "I"=1.029230775130462; "J"=1.028959929996554; "K"=1.031227397684115; "L"=1.052182979717557; "M"=1.056276083485348; "N"=1.030421171120503 }
    7 = @{ "B"=1.02; "C"=1.024172289829182; "D"=1.028552562660959; "E"=1.049383780435057; "F"=1.053460299732262; "I"=1.029177373959079; "J"=1.028759243118182; "K"=1.031054811601532; "L"=1.051833682185522; "M"=1.055900218072674; "N"=1.030220199243744 }
    8 = @{ "B"=1.02; "C"=1.022854740868141; "D"=1.027574472180263; "E"=1.047674165834534; "F"=1.051640245877543; "I"=1.028949464818814; "J"=1.02791785185789; "K"=1.030329743606959; "L"=1.050373454873257; "M"=1.05432876213717; "N"=1.029377613111354 }
    9 = @{ "B"=1.02; "C"=1.020527678471264; "D"=1.025842131748694; "E"=1.044664807429832; "F"=1.048435468803613; "I"=1.028530164685231; "J"=1.026425586997832; "K"=1.029038205797236; "L"=1.047799376176409; "M"=1.05155800706111; "N"=1.027883229064029 }
    10 = @{ "B"=1.02; "C"=1.018972656081374; "D"=1.024681333723657; "E"=1.04266056855094; "F"=1.046300397346572; "I"=1.028238754602909; "J"=1.025424265143034; "K"=1.028167866056484; "L"=1.046082537279574; "M"=1.049709593033122; "N"=1.026880485217345 }
    11 = @{ "B"=1.02; "C"=1.018298440015977; "D"=1.024177299146948; "E"=1.04179315089983; "F"=1.045376194484082; "I"=1.028109757505499; "J"=1.02498914248409; "K"=1.027788791985119; "L"=1.04533891030887; "M"=1.048908883443484; "N"=1.026444744634316 }
    12 = @{ "B"=1.02; "C"=1.018047872851104; "D"=1.023989867929502; "E"=1.041471015066514; "F"=1.045032945804039; "I"=1.028061419582641; "J"=1.024827286283051; "K"=1.027647654806617; "L"=1.045062657742417; "M"=1.048611411299885; "N"=1.026282658578906 }
    13 = @{ "B"=1.02; "C"=1.018101626405552; "D"=1.024030082076816; "E"=1.041540111553111; "F"=1.045106571989088; "I"=1.028071807354137; "J"=1.024862015504163; "K"=1.027677944242664; "L"=1.045121916546057; "M"=1.048675222524911; "N"=1.026317437119495 }
    14 = @{ "B"=1.02; "C"=1.018277730780604; "D"=1.02416181031366; "E"=1.041766521810286; "F"=1.045347820615033; "I"=1.028105770498651; "J"=1.02497576814037; "K"=1.027777132318438; "L"=1.045316075943654; "M"=1.048884295399429; "N"=1.026431351297493 }
    15 = @{ "B"=1.02; "C"=1.018386216807039; "D"=1.02424294456257; "E"=1.041906028727211; "F"=1.045496467246449; "I"=1.028126640312662; "J"=1.025045824104877; "K"=1.027838201347438; "L"=1.045435699052805; "M"=1.049013104998729; "N"=1.026501506749505 }
    16 = @{ "B"=1.02; "C"=1.019017382938848; "D"=1.024714755331326; "E"=1.042718145024396; "F"=1.046361739632057; "I"=1.028247256432709; "J"=1.025453110247595; "K"=1.028192977368462; "L"=1.046131884395621; "M"=1.049762726217283; "N"=1.026909371285262 }
    17 = @{ "B"=1.02; "C"=1.019413060234624; "D"=1.025010334727035; "E"=1.043227677365593; "F"=1.046904579421737; "I"=1.028322162336946; "J"=1.025708176553344; "K"=1.028414927007574; "L"=1.046568521097218; "M"=1.050232852410662; "N"=1.027164799814422 }
    18 = @{ "B"=1.02; "C"=1.019643767155401; "D"=1.025182605962949; "E"=1.043524920750869; "F"=1.047221237576776; "I"=1.028365582101543; "J"=1.025856803494892; "K"=1.028544173125907; "L"=1.046823182427325; "M"=1.050507037013043; "N"=1.027313637823275 }
    19 = @{ "B"=1.02; "C"=1.019722417871948; "D"=1.025241323059447; "E"=1.043626280244958; "F"=1.047329214822856; "I"=1.028380341047729; "J"=1.025907456192239; "K"=1.028588206523298; "L"=1.04691001179965; "M"=1.050600521523855; "N"=1.027364362453262 }
    20 = @{ "B"=1.02; "C"=1.019370616598776; "D"=1.024978635840457; "E"=1.043173005066546; "F"=1.046846334884424; "I"=1.028314153731731; "J"=1.025680825758409; "K"=1.028391135979374; "L"=1.046521676351303; "M"=1.050182415654186; "N"=1.027137410178221 }
    21 = @{ "B"=1.02; "C"=1.018225876118614; "D"=1.024123025435821; "E"=1.04169984795388; "F"=1.045276777790465; "I"=1.028095780864918; "J"=1.024942277238934; "K"=1.027747933072315; "L"=1.045258901866774; "M"=1.048822730132736; "N"=1.026397812835136 }
    22 = @{ "B"=1.02; "C"=1.017505360651482; "D"=1.023583852378534; "E"=1.0407739710918; "F"=1.044290172945404; "I"=1.027956035904602; "J"=1.024476578435754; "K"=1.027341603923712; "L"=1.044464733496713; "M"=1.047967534154477; "N"=1.025931452686238 }
    23 = @{ "B"=1.02; "C"=1.017887392992781; "D"=1.023869793496034; "E"=1.041264763259974; "F"=1.044813169364519; "I"=1.028030349070218; "J"=1.024723581702549; "K"=1.027557188870151; "L"=1.044885758163928; "M"=1.048420919814453; "N"=1.026178806726004 }
    24 = @{ "B"=1.02; "C"=1.019389795312966; "D"=1.024992959620506; "E"=1.043197708995257; "F"=1.046872652993646; "I"=1.028317773314766; "J"=1.025693184865242; "K"=1.028401886779988; "L"=1.046542843541232; "M"=1.05020520595173; "N"=1.027149786836403 }
    25 = @{ "B"=1.02; "C"=1.021129920137617; "D"=1.026291025768108; "E"=1.045442435016374; "F"=1.049263713723861; "I"=1.028640658615852; "J"=1.026812514586521; "K"=1.029373742017303; "L"=1.048464966049637; "M"=1.052274522972573; "N"=1.028270706134276 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
